$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 362, shifting existing rows 362..397 down to 363..398
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row 362 with the new record
$ws.Cells.Item(362, 1).Value = 10
$ws.Cells.Item(362, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(362, 3).Value = "La Araucanía"
$ws.Cells.Item(362, 4).Value = 44769
$ws.Cells.Item(362, 5).Value = 9
$ws.Cells.Item(362, 6).Value = 100112037
$ws.Cells.Item(362, 7).Value = "Cebollín"
$ws.Cells.Item(362, 8).Value = "Sin especificar"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 80
$ws.Cells.Item(362, 11).Value = 8000
$ws.Cells.Item(362, 12).Value = 8000
$ws.Cells.Item(362, 13).Value = 8000
$ws.Cells.Item(362, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(362, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(362, 16).Value = 667
$ws.Cells.Item(362, 17).Value = 12
$ws.Cells.Item(362, 18).Value = "Hortaliza"
